$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F98").Value = 13886963
$ws.Range("G98").Value = 11076399
$ws.Range("H98").Value = 904873
$ws.Range("M99").Value = 8586182
$ws.Range("M105").Value = 8358622
$ws.Range("H132").Value = 935378
$ws.Range("M132").Value = 14066925
$ws.Range("M133").Value = 14469377
$ws.Range("M134").Value = 14605013
$ws.Range("M135").Value = 14690449
$ws.Range("M136").Value = 14866909
$ws.Range("M137").Value = 14861427
$ws.Range("M138").Value = 14670818
$ws.Range("M139").Value = 14467044
$ws.Range("D140").Value = 823659
$ws.Range("M140").Value = 14131823
$ws.Range("N140").Value = 8935495
$ws.Range("M141").Value = 13894120
$ws.Range("N141").Value = 8782010
$ws.Range("M142").Value = 13750027
$ws.Range("N142").Value = 8776278
$ws.Range("M143").Value = 13621896
$ws.Range("N143").Value = 8836720
$ws.Range("N144").Value = 9042528
$ws.Range("N145").Value = 9395148
$ws.Range("N146").Value = 9716012
$ws.Range("N147").Value = 9916843
$ws.Range("N148").Value = 10176476
$ws.Range("N149").Value = 10402057
$ws.Range("N150").Value = 10586457
$ws.Range("N151").Value = 10834370
$ws.Range("J156").Value = 37.94
$ws.Range("D186").Value = 590414
$ws.Range("G186").Value = 9817610
$ws.Range("H186").Value = 1209130
$ws.Range("M186").Value = 13889919
$ws.Range("N186").Value = 8211219
$ws.Range("M187").Value = 13966152
$ws.Range("N187").Value = 8231210
$ws.Range("M188").Value = 14177385
$ws.Range("N188").Value = 8310996
$ws.Range("E189").Value = 204865
$ws.Range("M189").Value = 14260693
$ws.Range("N189").Value = 8315356
$ws.Range("O189").Value = 2521450
$ws.Range("M190").Value = 14423020
$ws.Range("N190").Value = 8370937
$ws.Range("O190").Value = 2540058
$ws.Range("D191").Value = 573547
$ws.Range("G191").Value = 8670436
$ws.Range("H191").Value = 1073058
$ws.Range("M191").Value = 14528175
$ws.Range("N191").Value = 8368458
$ws.Range("O191").Value = 2551089
$ws.Range("M192").Value = 14532950
$ws.Range("N192").Value = 8290032
$ws.Range("O192").Value = 2547711
$ws.Range("M193").Value = 14761473
$ws.Range("N193").Value = 8365124
$ws.Range("O193").Value = 2572359
$ws.Range("F194").Value = 14729057
$ws.Range("M194").Value = 14844154
$ws.Range("N194").Value = 8388011
$ws.Range("O194").Value = 2570159
$ws.Range("M195").Value = 14889085
$ws.Range("N195").Value = 8325483
$ws.Range("O195").Value = 2571048
$ws.Range("F196").Value = 13097796
$ws.Range("M196").Value = 14967767
$ws.Range("N196").Value = 8319979
$ws.Range("O196").Value = 2563122
$ws.Range("M197").Value = 15089911
$ws.Range("N197").Value = 8341443
$ws.Range("O197").Value = 2566577
$ws.Range("M198").Value = 15141506
$ws.Range("N198").Value = 8330307
$ws.Range("O198").Value = 2567378
$ws.Range("E199").Value = 214621
$ws.Range("F199").Value = 10187305
$ws.Range("M199").Value = 15309050
$ws.Range("N199").Value = 8356555
$ws.Range("E200").Value = 247050
$ws.Range("M200").Value = 15440660
$ws.Range("N200").Value = 8350574
$ws.Range("O200").Value = 2599975
$ws.Range("M201").Value = 15510434
$ws.Range("N201").Value = 8368317
$ws.Range("M202").Value = 15673744
$ws.Range("N202").Value = 8381874
$ws.Range("F203").Value = 9851475
$ws.Range("O211").Value = 2671977
$ws.Range("B232").Value = 1361144
$ws.Range("F232").Value = 11589074
$ws.Range("F625").Value = 7672046
$ws.Range("B626").Value = 1413965
$ws.Range("D626").Value = 528859
$ws.Range("E626").Value = 177858
$ws.Range("F626").Value = 9156916
$ws.Range("G626").Value = 7010215
$ws.Range("H626").Value = 2635461
$ws.Range("M626").Value = 41288564
$ws.Range("N626").Value = 7247067
$ws.Range("O626").Value = 4099366
$ws.Range("B627").Value = 905422
$ws.Range("D627").Value = 391934
$ws.Range("F627").Value = 7687898
$ws.Range("G627").Value = 6245143
$ws.Range("H627").Value = 2416130
$ws.Range("M627").Value = 38636537
$ws.Range("N627").Value = 6872875
$ws.Range("O627").Value = 3605418
$ws.Range("B628").Value = 892685
$ws.Range("D628").Value = 369746
$ws.Range("E628").Value = 168559
$ws.Range("F628").Value = 7553247
$ws.Range("G628").Value = 6442398
$ws.Range("H628").Value = 2542341
$ws.Range("I628").Value = 407.99
$ws.Range("L628").Value = 359.81
$ws.Range("M628").Value = 35832416
$ws.Range("N628").Value = 6392639
$ws.Range("O628").Value = 3165209
$ws.Range("B629").Value = 877570
$ws.Range("D629").Value = 304518
$ws.Range("E629").Value = 125434
$ws.Range("F629").Value = 6045368
$ws.Range("G629").Value = 4792485
$ws.Range("H629").Value = 1887960
$ws.Range("K629").Value = 17.4
$ws.Range("L629").Value = 363.47
$ws.Range("M629").Value = 33651738
$ws.Range("N629").Value = 5630884
$ws.Range("O629").Value = 2881911
$ws.Range("R629").Value = 323392
$ws.Range("B630").Value = 859776
$ws.Range("D630").Value = 312845
$ws.Range("F630").Value = 5997619
$ws.Range("G630").Value = 4677136
$ws.Range("H630").Value = 1815689
$ws.Range("K630").Value = 17.41
$ws.Range("L630").Value = 367.93
$ws.Range("M630").Value = 31549861
$ws.Range("N630").Value = 5191575
$ws.Range("O630").Value = 2637028
$ws.Range("P630").Value = 5258
$ws.Range("Q630").Value = 14105
$ws.Range("R630").Value = 320626
$ws.Range("B631").Value = 937470
$ws.Range("C631").Value = 1297.8
$ws.Range("D631").Value = 327913
$ws.Range("E631").Value = 106376
$ws.Range("F631").Value = 5526924
$ws.Range("G631").Value = 4384756
$ws.Range("H631").Value = 1693004
$ws.Range("I631").Value = 394.65
$ws.Range("J631").Value = 31.58
$ws.Range("L631").Value = 372.8
$ws.Range("M631").Value = 29255661
$ws.Range("N631").Value = 4860886
$ws.Range("O631").Value = 2409354
$ws.Range("P631").Value = 5250
$ws.Range("Q631").Value = 14116
$ws.Range("R631").Value = 318721
$ws.Range("A632").Value = 44773
$ws.Range("B632").Value = 1018530
$ws.Range("F632").Value = 6551392
$ws.Range("G632").Value = 4917713
$ws.Range("H632").Value = 1825095
$ws.Range("I632").Value = 380.45
$ws.Range("J632").Value = 30.93
$ws.Range("K632").Value = 16.54
$ws.Range("L632").Value = 377.26
$ws.Range("M632").Value = 27622875
$ws.Range("N632").Value = 4626533
$ws.Range("O632").Value = 2241562
$ws.Range("P632").Value = 5260
$ws.Range("Q632").Value = 14153
